$wb = $excel.ActiveWorkbook

# --- OFF sheet (Week "H" row, row 2): update target depth counts ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 215
$wsOff.Range("C2").Value = 157
$wsOff.Range("D2").Value = 69
$wsOff.Range("E2").Value = 38
$wsOff.Range("G2").Value = 4

# --- DEF sheet (Week "H" row, row 2): update target depth counts ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 169
$wsDef.Range("C2").Value = 113
$wsDef.Range("D2").Value = 39
$wsDef.Range("E2").Value = 15

$wb.Save()
